$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 20

# Read all existing data first (bottom-up doesn't matter since we read everything up front)
$segNames = @{}
$percAct  = @{}
$percArea = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $segNames[$r] = $ws.Cells.Item($r, 1).Value2
    $percAct[$r]  = $ws.Cells.Item($r, 2).Value2
    $percArea[$r] = $ws.Cells.Item($r, 3).Value2
}
$oldHeaderB1 = $ws.Cells.Item(1, 2).Value2
$oldHeaderC1 = $ws.Cells.Item(1, 3).Value2

# Copy the header style (from B1, which currently holds "PercActivations") to use for new D1
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Cells.Item(1, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Write new values: column layout becomes A=index, B=segments, C=PercActivations, D=PercSegmentAreas
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $percArea[$r]
    $ws.Cells.Item($r, 3).Value = $percAct[$r]
    $ws.Cells.Item($r, 2).Value = $segNames[$r]
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Header row
$ws.Cells.Item(1, 4).Value = $oldHeaderC1
$ws.Cells.Item(1, 3).Value = $oldHeaderB1
$ws.Cells.Item(1, 2).Value = "segments"

$ws.Range("A1").Select()
